# Apply the "new version with timestamp" update to the DaySale report.
#
# Summary of changes:
#  1) Row 26 ("برفان القصاص"): request-count 11:0 -> 10:0, sell price
#     30.0000 -> 60.0000, transactions 1:0 -> 2:0.
#  2) A new data row is inserted at row 39 for "مناديل بكر فاين"
#     (request 6:0, balance 0, price 15.00, sell price 15.0000,
#     transactions 1:0) with item # 33. The old "Total" row and the
#     footer row are pushed down from 39/40 to 40/41.
#  3) The Total row's P value is recomputed to reflect the new data.
#  4) The footer timestamp text is updated to the new generation time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Many columns in this sheet hold numeric-looking values ("30.0000",
    # "0", "10:0", ...) as literal TEXT (shared strings), not numbers.
    # Assigning such a string straight to .Value lets Excel auto-coerce
    # it to a real number/time, which both changes the stored type and
    # (worse) forks a brand new cell style. Temporarily flipping the
    # number format to "@" (text) for the assignment keeps the original
    # style id intact once the format is restored.
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = $fmt
}

# ---------------------------------------------------------------------
# 1) Update the existing "برفان القصاص" row (row 26).
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("H26") "10:0"
Set-TextValue $ws.Range("P26") "60.0000"
Set-TextValue $ws.Range("Q26") "2:0"

# ---------------------------------------------------------------------
# 2) Insert a new row at position 39 (pushes Total/footer rows down),
#    copy the formatting from the row above (row 38) so borders/fonts/
#    number formats match the rest of the data table, then fill it in.
# ---------------------------------------------------------------------
$ws.Range("A39:Q39").Insert(-4121)

$ws.Range("A38:Q38").Copy()
$ws.Range("A39:Q39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(39).RowHeight = 25.5

$ws.Range("A39").Value = 33
$ws.Range("C39").Value = "مناديل بكر فاين"
Set-TextValue $ws.Range("H39") "6:0"
Set-TextValue $ws.Range("L39") "0"
Set-TextValue $ws.Range("N39") "15.00"
Set-TextValue $ws.Range("P39") "15.0000"
Set-TextValue $ws.Range("Q39") "1:0"

$ws.Range("A39:B39").Merge()
$ws.Range("C39:G39").Merge()
$ws.Range("H39:K39").Merge()
$ws.Range("L39:M39").Merge()
$ws.Range("N39:O39").Merge()

# ---------------------------------------------------------------------
# 3) The Total row (now row 40) - recompute the grand total of the
#    "sell price" column to account for the changed/added rows.
#    (30 -> 60 for row 26, plus the new 15.0000 row = +45 overall)
# ---------------------------------------------------------------------
$ws.Rows.Item(40).RowHeight = 24.75
$ws.Range("P40").Value = 1055.5450000000001

# ---------------------------------------------------------------------
# 4) Footer row (now row 41) - refresh the generated-on timestamp text.
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("A41") "Thursday, 24 July, 2025 2:07 PM"
